$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions refresh).
# Columns B (Coin) and C (Link) are plain text; D (Price) and E (Volume)
# are forced to text format ("@") before assignment so Excel does not
# reinterpret numeric-looking strings (e.g. "6.000", "17.60", "28.027.82")
# as numbers/dates and strip significant formatting.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.027.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.46%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.864.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.98%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.35%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5122'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.57%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3855'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.41%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08343'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.21%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.90%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.36'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.174'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.45'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.37'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.59%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.253'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.30%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.10%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.61%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.79'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.80%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06637'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.60'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.000'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.71%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.065.37'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.45%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.98'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.29%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.238'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.52%  '

$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.075.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.78%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.459'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.59%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.82'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.00%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.49'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.46%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.84'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.50%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1062'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.10%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.026'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.62%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.800'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.76%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.599'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.31%  '

$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.433'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.82%  '

$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02417'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.92%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06508'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.85%  '

$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2174'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.39%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.197'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.85%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6463'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.70%  '

$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.974'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.99%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.209'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.36%  '

$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.27'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.61%  '

$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6060'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.29%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.01'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.93%  '

$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.286'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.74%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.674'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.38%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.997'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.34%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.219'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.39%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.49'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.43%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.08'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.26%  '
